# chore: update Sheets via scheduled runner
# Refreshes cached market-board price/profit figures across the per-job
# leve-profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8211.75
$ws.Range("I51").Value = 4499.6665
$ws.Range("J51").Value = 10439
$ws.Range("K51").Value = 4499.6665
$ws.Range("L51").Value = 10439
$ws.Range("M51").Value = -4015.6665
$ws.Range("N51").Value = -11407

$ws.Range("H123").Value = 57215
$ws.Range("J123").Value = 57215
$ws.Range("L123").Value = 57215
$ws.Range("N123").Value = -67015

$ws.Range("I125").Value = 29412242
$ws.Range("J125").Value = 699.1429000000001
$ws.Range("K125").Value = 264710178
$ws.Range("L125").Value = 6292.2861
$ws.Range("M125").Value = -264707718
$ws.Range("N125").Value = -11212.2861

$ws.Range("H132").Value = 1862.3489
$ws.Range("I132").Value = 1787.6428
$ws.Range("K132").Value = 5362.928400000001
$ws.Range("M132").Value = -2832.928400000001

$ws.Range("H133").Value = 100780
$ws.Range("J133").Value = 100780
$ws.Range("L133").Value = 100780
$ws.Range("N133").Value = -110900

$ws.Range("H137").Value = 2549.3696
$ws.Range("I137").Value = 2894.125
$ws.Range("J137").Value = 2365.5
$ws.Range("K137").Value = 8682.375
$ws.Range("L137").Value = 7096.5
$ws.Range("M137").Value = -6132.375
$ws.Range("N137").Value = -12196.5

$ws.Range("H138").Value = 1646034.2
$ws.Range("J138").Value = 2007344.1
$ws.Range("L138").Value = 6022032.300000001
$ws.Range("N138").Value = -6032312.300000001

$ws.Range("H141").Value = 4188
$ws.Range("I141").Value = 2982.5
$ws.Range("J141").Value = 6599
$ws.Range("K141").Value = 8947.5
$ws.Range("L141").Value = 19797
$ws.Range("M141").Value = -3767.5
$ws.Range("N141").Value = -30157

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2035.2646
$ws.Range("I2").Value = 1664.8572
$ws.Range("K2").Value = 1664.8572
$ws.Range("M2").Value = -1551.8572

$ws.Range("H32").Value = 3579.275
$ws.Range("I32").Value = 3022.514
$ws.Range("K32").Value = 3022.514
$ws.Range("M32").Value = -2735.514

$ws.Range("H116").Value = 2035.2646
$ws.Range("I116").Value = 1664.8572
$ws.Range("K116").Value = 1664.8572
$ws.Range("M116").Value = 629.1428000000001

$ws.Range("H122").Value = 10632.828
$ws.Range("I122").Value = 13359.73
$ws.Range("J122").Value = 2755.111
$ws.Range("K122").Value = 40079.19
$ws.Range("L122").Value = 8265.332999999999
$ws.Range("M122").Value = -37629.19
$ws.Range("N122").Value = -13165.333

$ws.Range("H132").Value = 1389664.6
$ws.Range("I132").Value = 2107154.2
$ws.Range("K132").Value = 6321462.600000001
$ws.Range("M132").Value = -6318932.600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2035.2646
$ws.Range("I3").Value = 1664.8572
$ws.Range("K3").Value = 1664.8572
$ws.Range("M3").Value = -1550.8572

$ws.Range("H20").Value = 5557310
$ws.Range("I20").Value = 9805492
$ws.Range("J20").Value = 1995.5385
$ws.Range("K20").Value = 9805492
$ws.Range("L20").Value = 1995.5385
$ws.Range("M20").Value = -9805245
$ws.Range("N20").Value = -2489.5385

$ws.Range("H105").Value = 1113926.2
$ws.Range("I105").Value = 1819686.9
$ws.Range("J105").Value = 4873.857
$ws.Range("K105").Value = 1819686.9
$ws.Range("L105").Value = 4873.857
$ws.Range("M105").Value = -1817939.9
$ws.Range("N105").Value = -8367.857

$ws.Range("H107").Value = 102276824
$ws.Range("I107").Value = 160717710
$ws.Range("J107").Value = 5257.25
$ws.Range("K107").Value = 160717710
$ws.Range("L107").Value = 5257.25
$ws.Range("M107").Value = -160715790
$ws.Range("N107").Value = -9097.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6349.056
$ws.Range("I31").Value = 4451.643
$ws.Range("J31").Value = 7584.5815
$ws.Range("K31").Value = 4451.643
$ws.Range("L31").Value = 7584.5815
$ws.Range("M31").Value = -4156.643
$ws.Range("N31").Value = -8174.5815

$ws.Range("H34").Value = 6349.056
$ws.Range("I34").Value = 4451.643
$ws.Range("J34").Value = 7584.5815
$ws.Range("K34").Value = 4451.643
$ws.Range("L34").Value = 7584.5815
$ws.Range("M34").Value = -4249.643
$ws.Range("N34").Value = -7988.5815

$ws.Range("H132").Value = 5579.273
$ws.Range("I132").Value = 3033.2
$ws.Range("J132").Value = 9496.308000000001
$ws.Range("K132").Value = 9099.599999999999
$ws.Range("L132").Value = 28488.924
$ws.Range("M132").Value = -6569.599999999999
$ws.Range("N132").Value = -33548.924

$ws.Range("H134").Value = 5324.9565
$ws.Range("I134").Value = 2250.4285
$ws.Range("K134").Value = 6751.2855
$ws.Range("M134").Value = -4216.2855

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H140").Value = 60000
$ws.Range("J140").Value = 60000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 666666700
$ws.Range("J32").Value = 1000000000
$ws.Range("L32").Value = 3000000000
$ws.Range("N32").Value = -3000000566

$ws.Range("H68").Value = 3339.027
$ws.Range("I68").Value = 2528.4285
$ws.Range("J68").Value = 3528.1667
$ws.Range("K68").Value = 7585.2855
$ws.Range("L68").Value = 10584.5001
$ws.Range("M68").Value = -6774.2855
$ws.Range("N68").Value = -12206.5001

$ws.Range("H71").Value = 3339.027
$ws.Range("I71").Value = 2528.4285
$ws.Range("J71").Value = 3528.1667
$ws.Range("K71").Value = 22755.8565
$ws.Range("L71").Value = 31753.5003
$ws.Range("M71").Value = -18699.8565
$ws.Range("N71").Value = -39865.5003

$ws.Range("H86").Value = 425
$ws.Range("I86").Value = 425
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1275
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -89
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 425
$ws.Range("I89").Value = 425
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 3825
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 2103
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2858.0278
$ws.Range("I113").Value = 1844.5385
$ws.Range("J113").Value = 5493.1
$ws.Range("K113").Value = 1844.5385
$ws.Range("L113").Value = 5493.1
$ws.Range("M113").Value = 325.4614999999999
$ws.Range("N113").Value = -9833.1

$ws.Range("H132").Value = 6050
$ws.Range("I132").Value = 3191.182
$ws.Range("J132").Value = 10542.429
$ws.Range("K132").Value = 9573.545999999998
$ws.Range("L132").Value = 31627.287
$ws.Range("M132").Value = -7043.545999999998
$ws.Range("N132").Value = -36687.287

$ws.Range("H133").Value = 98593.336
$ws.Range("J133").Value = 98593.336
$ws.Range("L133").Value = 98593.336
$ws.Range("N133").Value = -108713.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 16136633
$ws.Range("I132").Value = 35717908
$ws.Range("J132").Value = 10876.059
$ws.Range("K132").Value = 107153724
$ws.Range("L132").Value = 32628.177
$ws.Range("M132").Value = -107151194
$ws.Range("N132").Value = -37688.177

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 56959
$ws.Range("J121").Value = 56959
$ws.Range("L121").Value = 56959
$ws.Range("N121").Value = -60453

$ws.Range("H132").Value = 26341.5
$ws.Range("I132").Value = 10525.237
$ws.Range("K132").Value = 31575.711
$ws.Range("M132").Value = -29045.711
